# Sync attendance_reports: normalize the "Recorded By" (column G) audit
# trail entries. Several rows list the automated "System" editor together
# with a human editor (and a couple of rows list two human editors); this
# re-orders each such comma-separated list to the canonical order used
# across the rest of the sheet. Only cells whose current text exactly
# matches one of the known stale orderings are touched - everything else
# (single-author cells, already-canonical cells) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "System, backup@backdoor.com, system" = "system, backup@backdoor.com, System";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, backup@backdoor.com"         = "backup@backdoor.com, System";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G = "Recorded By" (row 1 is the header row, so data starts at row 2).
for ($r = [Math]::Max(2, $firstRow); $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($null -eq $val) { continue }

    if ($replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}
